$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1948051948051948
$ws.Range("C2").Value = 0.5281385281385281
$ws.Range("J2").Value = 0.01731601731601732
$ws.Range("P2").Value = 0.1385281385281385
$ws.Range("S2").Value = 0.1212121212121212
$ws.Range("B3").Value = 0.03125
$ws.Range("C3").Value = 0.0546875
$ws.Range("J3").Value = 0.046875
$ws.Range("P3").Value = 0.671875
$ws.Range("S3").Value = 0.1953125
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.5789473684210527
$ws.Range("S4").Value = 0.3421052631578947
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.02314814814814815
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("F6").Value = 0.08796296296296297
$ws.Range("J6").Value = 0.1990740740740741
$ws.Range("O6").Value = 0.01388888888888889
$ws.Range("Q6").Value = 0.2268518518518519
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.3796296296296297
$ws.Range("B7").Value = 0.08
$ws.Range("D7").Value = 0.02285714285714286
$ws.Range("F7").Value = 0.06857142857142857
$ws.Range("J7").Value = 0.1714285714285714
$ws.Range("O7").Value = 0.01142857142857143
$ws.Range("Q7").Value = 0.1771428571428571
$ws.Range("R7").Value = 0.09142857142857143
$ws.Range("S7").Value = 0.3771428571428572
$ws.Range("B8").Value = 0.07526881720430108
$ws.Range("D8").Value = 0.01720430107526882
$ws.Range("E8").Value = 0.002150537634408602
$ws.Range("F8").Value = 0.06236559139784946
$ws.Range("J8").Value = 0.09462365591397849
$ws.Range("O8").Value = 0.02365591397849462
$ws.Range("Q8").Value = 0.2150537634408602
$ws.Range("R8").Value = 0.07741935483870968
$ws.Range("S8").Value = 0.432258064516129
$ws.Range("B9").Value = 0.08648648648648649
$ws.Range("D9").Value = 0.02162162162162162
$ws.Range("F9").Value = 0.05945945945945946
$ws.Range("J9").Value = 0.08108108108108109
$ws.Range("O9").Value = 0.01621621621621622
$ws.Range("Q9").Value = 0.2054054054054054
$ws.Range("R9").Value = 0.0918918918918919
$ws.Range("S9").Value = 0.4378378378378379
$ws.Range("B10").Value = 0.09821428571428571
$ws.Range("D10").Value = 0.01964285714285714
$ws.Range("E10").Value = 0.002678571428571429
$ws.Range("F10").Value = 0.08214285714285714
$ws.Range("J10").Value = 0.1276785714285714
$ws.Range("O10").Value = 0.00625
$ws.Range("Q10").Value = 0.2223214285714286
$ws.Range("R10").Value = 0.05446428571428572
$ws.Range("S10").Value = 0.3866071428571429
$ws.Range("G11").Value = 0.1654676258992806
$ws.Range("J11").Value = 0.09352517985611511
$ws.Range("K11").Value = 0.2014388489208633
$ws.Range("L11").Value = 0.5287769784172662
$ws.Range("S11").Value = 0.01079136690647482
$ws.Range("G12").Value = 0.7105263157894737
$ws.Range("J12").Value = 0.1973684210526316
$ws.Range("K12").Value = 0.006578947368421052
$ws.Range("L12").Value = 0.03289473684210526
$ws.Range("S12").Value = 0.05263157894736842
$ws.Range("G13").Value = 0.6444444444444445
$ws.Range("J13").Value = 0.3111111111111111
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.005181347150259068
$ws.Range("H15").Value = 0.1761658031088083
$ws.Range("I15").Value = 0.07253886010362694
$ws.Range("J15").Value = 0.3419689119170984
$ws.Range("K15").Value = 0.07253886010362694
$ws.Range("M15").Value = 0.0155440414507772
$ws.Range("O15").Value = 0.07253886010362694
$ws.Range("S15").Value = 0.2435233160621762
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.08571428571428572
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.09285714285714286
$ws.Range("M16").Value = 0.04285714285714286
$ws.Range("N16").Value = 0.007142857142857143
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.09285714285714286
$ws.Range("F17").Value = 0.01521739130434783
$ws.Range("H17").Value = 0.1869565217391304
$ws.Range("I17").Value = 0.08695652173913043
$ws.Range("J17").Value = 0.3891304347826087
$ws.Range("K17").Value = 0.1065217391304348
$ws.Range("M17").Value = 0.02826086956521739
$ws.Range("O17").Value = 0.04782608695652174
$ws.Range("S17").Value = 0.1391304347826087
$ws.Range("F18").Value = 0.007092198581560284
$ws.Range("H18").Value = 0.198581560283688
$ws.Range("I18").Value = 0.0851063829787234
$ws.Range("J18").Value = 0.4468085106382979
$ws.Range("K18").Value = 0.1134751773049645
$ws.Range("M18").Value = 0.02836879432624113
$ws.Range("O18").Value = 0.03546099290780142
$ws.Range("S18").Value = 0.0851063829787234
$ws.Range("F19").Value = 0.01237623762376238
$ws.Range("H19").Value = 0.2409240924092409
$ws.Range("I19").Value = 0.08828382838283828
$ws.Range("J19").Value = 0.3432343234323432
$ws.Range("K19").Value = 0.1014851485148515
$ws.Range("M19").Value = 0.01567656765676568
$ws.Range("O19").Value = 0.07673267326732673
$ws.Range("S19").Value = 0.1212871287128713
